# Apply update described by commit:
# "Atualizacao de bases das ligas, do dia: 29-05-2024 as 22:54"
# Rewrites match rows 298-306 on sheet "Poland Ekstraklasa" with refreshed
# odds/results data (rows reshuffled among themselves; ids, teams, scores,
# and odds columns updated to the new snapshot).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland Ekstraklasa")

$ws.Range("B298").Value = 7083188
$ws.Range("E298").Value = 'Legia Warsaw'
$ws.Range("F298").Value = 'Zaglebie Lubin'
$ws.Range("G298").Value = 2
$ws.Range("H298").Value = 1
$ws.Range("I298").Value = 2
$ws.Range("K298").Value = 'H'
$ws.Range("L298").Value = 1.5
$ws.Range("M298").Value = 4
$ws.Range("N298").Value = 5.5
$ws.Range("O298").Value = 1.6
$ws.Range("P298").Value = 4.1
$ws.Range("Q298").Value = 4.333
$ws.Range("R298").Value = -0.75
$ws.Range("S298").Value = 1.825
$ws.Range("T298").Value = 2.025
$ws.Range("U298").Value = 3
$ws.Range("V298").Value = 1.875
$ws.Range("W298").Value = 1.975
$ws.Range("X298").Value = 0.6000000000000001
$ws.Range("Z298").Value = -1
$ws.Range("AA298").Value = 0.4125
$ws.Range("AB298").Value = -0.5
$ws.Range("AC298").Value = 0
$ws.Range("AD298").Value = 0
$ws.Range("B299").Value = 7093820
$ws.Range("E299").Value = 'Ruch Chorzow'
$ws.Range("F299").Value = 'Cracovia Krakow'
$ws.Range("G299").Value = 2
$ws.Range("I299").Value = 1
$ws.Range("L299").Value = 2.5
$ws.Range("M299").Value = 3.4
$ws.Range("N299").Value = 2.5
$ws.Range("O299").Value = 2.6
$ws.Range("P299").Value = 3.6
$ws.Range("Q299").Value = 2.3
$ws.Range("R299").Value = 0
$ws.Range("S299").Value = 2.025
$ws.Range("T299").Value = 1.825
$ws.Range("U299").Value = 3
$ws.Range("X299").Value = 1.6
$ws.Range("AA299").Value = 1.025
$ws.Range("AB299").Value = -1
$ws.Range("B300").Value = 7074364
$ws.Range("E300").Value = 'Rakow Czestochowa'
$ws.Range("F300").Value = 'Slask Wroclaw'
$ws.Range("G300").Value = 1
$ws.Range("H300").Value = 2
$ws.Range("K300").Value = 'A'
$ws.Range("M300").Value = 3.6
$ws.Range("N300").Value = 2.4
$ws.Range("O300").Value = 2.15
$ws.Range("Q300").Value = 2.875
$ws.Range("R300").Value = -0.25
$ws.Range("S300").Value = 1.95
$ws.Range("T300").Value = 1.9
$ws.Range("U300").Value = 2.5
$ws.Range("V300").Value = 1.875
$ws.Range("W300").Value = 1.975
$ws.Range("X300").Value = -1
$ws.Range("Z300").Value = 1.875
$ws.Range("AA300").Value = -1
$ws.Range("AB300").Value = 0.8999999999999999
$ws.Range("AC300").Value = 0.875
$ws.Range("AD300").Value = -1
$ws.Range("B301").Value = 7083189
$ws.Range("E301").Value = 'Pogon Szczecin'
$ws.Range("F301").Value = 'Gornik Zabrze'
$ws.Range("G301").Value = 1
$ws.Range("H301").Value = 0
$ws.Range("I301").Value = 0
$ws.Range("L301").Value = 1.727
$ws.Range("M301").Value = 4
$ws.Range("N301").Value = 3.75
$ws.Range("O301").Value = 1.55
$ws.Range("P301").Value = 4.333
$ws.Range("Q301").Value = 4.5
$ws.Range("R301").Value = -1
$ws.Range("S301").Value = 1.925
$ws.Range("T301").Value = 1.925
$ws.Range("U301").Value = 3.5
$ws.Range("V301").Value = 2.025
$ws.Range("W301").Value = 1.825
$ws.Range("X301").Value = 0.55
$ws.Range("AA301").Value = 0
$ws.Range("AB301").Value = 0
$ws.Range("AC301").Value = -1
$ws.Range("AD301").Value = 0.825
$ws.Range("B302").Value = 7041338
$ws.Range("E302").Value = 'Jagiellonia Bialystok'
$ws.Range("F302").Value = 'Warta Poznan'
$ws.Range("G302").Value = 3
$ws.Range("H302").Value = 0
$ws.Range("I302").Value = 3
$ws.Range("L302").Value = 1.444
$ws.Range("M302").Value = 4.75
$ws.Range("N302").Value = 5.25
$ws.Range("O302").Value = 1.4
$ws.Range("P302").Value = 4.75
$ws.Range("Q302").Value = 5.75
$ws.Range("R302").Value = -1.25
$ws.Range("S302").Value = 1.9
$ws.Range("T302").Value = 1.95
$ws.Range("V302").Value = 1.925
$ws.Range("W302").Value = 1.925
$ws.Range("X302").Value = 0.3999999999999999
$ws.Range("AA302").Value = 0.8999999999999999
$ws.Range("AB302").Value = -1
$ws.Range("B303").Value = 7093821
$ws.Range("E303").Value = 'LKS Lodz'
$ws.Range("F303").Value = 'Stal Mielec'
$ws.Range("G303").Value = 3
$ws.Range("H303").Value = 2
$ws.Range("I303").Value = 3
$ws.Range("K303").Value = 'H'
$ws.Range("L303").Value = 2.5
$ws.Range("M303").Value = 3.4
$ws.Range("N303").Value = 2.5
$ws.Range("O303").Value = 2.2
$ws.Range("P303").Value = 3.5
$ws.Range("Q303").Value = 2.8
$ws.Range("S303").Value = 2.025
$ws.Range("T303").Value = 1.825
$ws.Range("U303").Value = 3
$ws.Range("V303").Value = 2
$ws.Range("W303").Value = 1.85
$ws.Range("X303").Value = 1.2
$ws.Range("Z303").Value = -1
$ws.Range("AA303").Value = 1.025
$ws.Range("AB303").Value = -1
$ws.Range("AC303").Value = 1
$ws.Range("B304").Value = 7083187
$ws.Range("E304").Value = 'Lech Poznan'
$ws.Range("F304").Value = 'Korona Kielce'
$ws.Range("H304").Value = 2
$ws.Range("I304").Value = 1
$ws.Range("K304").Value = 'A'
$ws.Range("L304").Value = 1.8
$ws.Range("M304").Value = 3.8
$ws.Range("N304").Value = 3.6
$ws.Range("O304").Value = 2.1
$ws.Range("P304").Value = 3.7
$ws.Range("Q304").Value = 2.9
$ws.Range("R304").Value = -0.25
$ws.Range("S304").Value = 1.9
$ws.Range("T304").Value = 1.95
$ws.Range("U304").Value = 2.75
$ws.Range("V304").Value = 1.925
$ws.Range("W304").Value = 1.925
$ws.Range("X304").Value = -1
$ws.Range("Z304").Value = 1.9
$ws.Range("AA304").Value = -1
$ws.Range("AB304").Value = 0.95
$ws.Range("AC304").Value = 0.4625
$ws.Range("AD304").Value = -0.5
$ws.Range("B305").Value = 7090293
$ws.Range("E305").Value = 'Radomiak Radom'
$ws.Range("F305").Value = 'Widzew Lodz'
$ws.Range("G305").Value = 1
$ws.Range("H305").Value = 3
$ws.Range("I305").Value = 1
$ws.Range("K305").Value = 'A'
$ws.Range("L305").Value = 2.2
$ws.Range("M305").Value = 3.1
$ws.Range("N305").Value = 3.1
$ws.Range("O305").Value = 2.15
$ws.Range("P305").Value = 3.2
$ws.Range("Q305").Value = 3.1
$ws.Range("R305").Value = -0.25
$ws.Range("S305").Value = 1.925
$ws.Range("T305").Value = 1.925
$ws.Range("U305").Value = 2.75
$ws.Range("V305").Value = 1.9
$ws.Range("W305").Value = 1.95
$ws.Range("X305").Value = -1
$ws.Range("Z305").Value = 2.1
$ws.Range("AA305").Value = -1
$ws.Range("AB305").Value = 0.925
$ws.Range("AC305").Value = 0.8999999999999999
$ws.Range("B306").Value = 7088350
$ws.Range("E306").Value = 'Puszcza Niepolomice'
$ws.Range("F306").Value = 'Piast Gliwice'
$ws.Range("H306").Value = 0
$ws.Range("I306").Value = 0
$ws.Range("K306").Value = 'H'
$ws.Range("L306").Value = 3
$ws.Range("M306").Value = 3.1
$ws.Range("N306").Value = 2.3
$ws.Range("O306").Value = 2.7
$ws.Range("P306").Value = 3
$ws.Range("Q306").Value = 2.625
$ws.Range("R306").Value = 0
$ws.Range("S306").Value = 1.975
$ws.Range("T306").Value = 1.875
$ws.Range("U306").Value = 2.25
$ws.Range("V306").Value = 2.025
$ws.Range("W306").Value = 1.825
$ws.Range("X306").Value = 1.7
$ws.Range("Z306").Value = -1
$ws.Range("AA306").Value = 0.9750000000000001
$ws.Range("AB306").Value = -1
$ws.Range("AC306").Value = -1
$ws.Range("AD306").Value = 0.825
